$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Notes column text: "notes" -> "note"
$ws.Range("C2").Value = "{d.records[i].note}"
$ws.Range("C3").Value = "{d.records[i+1].note}"

# Update selection to match the new active cell
$ws.Range("E8").Select()
